# Align FW10 to FW9 settings
# The content-type / document-management schema stored in the package's
# customXml part (rId4 -> customXml/item4.xml) is refreshed to match the
# FW9 field layout: bump the schema/content-type version stamps, re-order
# a few ns3 field refs to their FW9 positions, add the
# MediaServiceSearchProperties field, and roll the TaxCatchAll taxonomy
# list GUID to the FW9 list.
$d = $word.ActiveDocument

$newSchemaXml = '<?xml version="1.0" encoding="utf-8"?><ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100B3416DCB34DD7C4B92E65D12CB7FCD38" ma:contentTypeVersion="16" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="bac8be2826a85e5a9dd452e9d6c31841"><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="79cd7128-7692-472f-9aff-5e5a6b9f80d3" xmlns:ns3="4c813fe8-1347-4959-8219-8a1b0e95da11" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="bef6fd18dd0aaa5b9dabdafd5fa7b70b" ns2:_="" ns3:_=""><xsd:import namespace="79cd7128-7692-472f-9aff-5e5a6b9f80d3"/><xsd:import namespace="4c813fe8-1347-4959-8219-8a1b0e95da11"/><xsd:element name="properties"><xsd:complexType><xsd:sequence><xsd:element name="documentManagement"><xsd:complexType><xsd:all><xsd:element ref="ns2:OX2_Project_Id" minOccurs="0"/><xsd:element ref="ns2:SharedWithUsers" minOccurs="0"/><xsd:element ref="ns2:SharedWithDetails" minOccurs="0"/><xsd:element ref="ns3:MediaServiceMetadata" minOccurs="0"/><xsd:element ref="ns3:MediaServiceFastMetadata" minOccurs="0"/><xsd:element ref="ns3:MediaServiceObjectDetectorVersions" minOccurs="0"/><xsd:element ref="ns3:lcf76f155ced4ddcb4097134ff3c332f" minOccurs="0"/><xsd:element ref="ns2:TaxCatchAll" minOccurs="0"/><xsd:element ref="ns3:MediaServiceGenerationTime" minOccurs="0"/><xsd:element ref="ns3:MediaServiceEventHashCode" minOccurs="0"/><xsd:element ref="ns3:MediaServiceOCR" minOccurs="0"/><xsd:element ref="ns3:MediaServiceDateTaken" minOccurs="0"/><xsd:element ref="ns3:MediaLengthInSeconds" minOccurs="0"/><xsd:element ref="ns3:MediaServiceLocation" minOccurs="0"/><xsd:element ref="ns3:MediaServiceSearchProperties" minOccurs="0"/></xsd:all></xsd:complexType></xsd:element></xsd:sequence></xsd:complexType></xsd:element></xsd:schema><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="79cd7128-7692-472f-9aff-5e5a6b9f80d3" elementFormDefault="qualified"><xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/><xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/><xsd:element name="OX2_Project_Id" ma:index="8" nillable="true" ma:displayName="Project Number" ma:internalName="OX2_Project_Id"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="SharedWithUsers" ma:index="9" nillable="true" ma:displayName="Shared With" ma:internalName="SharedWithUsers" ma:readOnly="true"><xsd:complexType><xsd:complexContent><xsd:extension base="dms:UserMulti"><xsd:sequence><xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded"><xsd:complexType><xsd:sequence><xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/><xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/><xsd:element name="AccountType" type="xsd:string" minOccurs="0"/></xsd:sequence></xsd:complexType></xsd:element></xsd:sequence></xsd:extension></xsd:complexContent></xsd:complexType></xsd:element><xsd:element name="SharedWithDetails" ma:index="10" nillable="true" ma:displayName="Shared With Details" ma:internalName="SharedWithDetails" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"><xsd:maxLength value="255"/></xsd:restriction></xsd:simpleType></xsd:element><xsd:element name="TaxCatchAll" ma:index="16" nillable="true" ma:displayName="Taxonomy Catch All Column" ma:hidden="true" ma:list="{6ec4a4bc-2685-4191-8e5b-741c43db70a3}" ma:internalName="TaxCatchAll" ma:showField="CatchAllData" ma:web="79cd7128-7692-472f-9aff-5e5a6b9f80d3"><xsd:complexType><xsd:complexContent><xsd:extension base="dms:MultiChoiceLookup"><xsd:sequence><xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/></xsd:sequence></xsd:extension></xsd:complexContent></xsd:complexType></xsd:element></xsd:schema><xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="4c813fe8-1347-4959-8219-8a1b0e95da11" elementFormDefault="qualified"><xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/><xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/><xsd:element name="MediaServiceMetadata" ma:index="11" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceFastMetadata" ma:index="12" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceObjectDetectorVersions" ma:index="13" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="lcf76f155ced4ddcb4097134ff3c332f" ma:index="15" nillable="true" ma:taxonomy="true" ma:internalName="lcf76f155ced4ddcb4097134ff3c332f" ma:taxonomyFieldName="MediaServiceImageTags" ma:displayName="Image Tags" ma:readOnly="false" ma:fieldId="{5cf76f15-5ced-4ddc-b409-7134ff3c332f}" ma:taxonomyMulti="true" ma:sspId="d1a101f0-c96a-436f-92e5-2eae4ddf4f0a" ma:termSetId="09814cd3-568e-fe90-9814-8d621ff8fb84" ma:anchorId="fba54fb3-c3e1-fe81-a776-ca4b69148c4d" ma:open="true" ma:isKeyword="false"><xsd:complexType><xsd:sequence><xsd:element ref="pc:Terms" minOccurs="0" maxOccurs="1"/></xsd:sequence></xsd:complexType></xsd:element><xsd:element name="MediaServiceGenerationTime" ma:index="17" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceEventHashCode" ma:index="18" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceOCR" ma:index="19" nillable="true" ma:displayName="Extracted Text" ma:internalName="MediaServiceOCR" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"><xsd:maxLength value="255"/></xsd:restriction></xsd:simpleType></xsd:element><xsd:element name="MediaServiceDateTaken" ma:index="20" nillable="true" ma:displayName="MediaServiceDateTaken" ma:description="" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaLengthInSeconds" ma:index="21" nillable="true" ma:displayName="MediaLengthInSeconds" ma:hidden="true" ma:internalName="MediaLengthInSeconds" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Unknown"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceLocation" ma:index="22" nillable="true" ma:displayName="Location" ma:description="" ma:indexed="true" ma:internalName="MediaServiceLocation" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Text"/></xsd:simpleType></xsd:element><xsd:element name="MediaServiceSearchProperties" ma:index="23" nillable="true" ma:displayName="MediaServiceSearchProperties" ma:hidden="true" ma:internalName="MediaServiceSearchProperties" ma:readOnly="true"><xsd:simpleType><xsd:restriction base="dms:Note"/></xsd:simpleType></xsd:element></xsd:schema><xsd:schema xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd" targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all"><xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/><xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/><xsd:element name="coreProperties" type="CT_coreProperties"/><xsd:complexType name="CT_coreProperties"><xsd:all><xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/><xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/><xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/><xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/><xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/><xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/><xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/><xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/><xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string"><xsd:annotation><xsd:documentation>This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.</xsd:documentation></xsd:annotation></xsd:element><xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/><xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/><xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/></xsd:all></xsd:complexType></xsd:schema><xs:schema xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema" targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified"><xs:element name="Person"><xs:complexType><xs:sequence><xs:element ref="pc:DisplayName" minOccurs="0"/><xs:element ref="pc:AccountId" minOccurs="0"/><xs:element ref="pc:AccountType" minOccurs="0"/></xs:sequence></xs:complexType></xs:element><xs:element name="DisplayName" type="xs:string"/><xs:element name="AccountId" type="xs:string"/><xs:element name="AccountType" type="xs:string"/><xs:element name="BDCAssociatedEntity"><xs:complexType><xs:sequence><xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"/></xs:sequence><xs:attribute ref="pc:EntityNamespace"/><xs:attribute ref="pc:EntityName"/><xs:attribute ref="pc:SystemInstanceName"/><xs:attribute ref="pc:AssociationName"/></xs:complexType></xs:element><xs:attribute name="EntityNamespace" type="xs:string"/><xs:attribute name="EntityName" type="xs:string"/><xs:attribute name="SystemInstanceName" type="xs:string"/><xs:attribute name="AssociationName" type="xs:string"/><xs:element name="BDCEntity"><xs:complexType><xs:sequence><xs:element ref="pc:EntityDisplayName" minOccurs="0"/><xs:element ref="pc:EntityInstanceReference" minOccurs="0"/><xs:element ref="pc:EntityId1" minOccurs="0"/><xs:element ref="pc:EntityId2" minOccurs="0"/><xs:element ref="pc:EntityId3" minOccurs="0"/><xs:element ref="pc:EntityId4" minOccurs="0"/><xs:element ref="pc:EntityId5" minOccurs="0"/></xs:sequence></xs:complexType></xs:element><xs:element name="EntityDisplayName" type="xs:string"/><xs:element name="EntityInstanceReference" type="xs:string"/><xs:element name="EntityId1" type="xs:string"/><xs:element name="EntityId2" type="xs:string"/><xs:element name="EntityId3" type="xs:string"/><xs:element name="EntityId4" type="xs:string"/><xs:element name="EntityId5" type="xs:string"/><xs:element name="Terms"><xs:complexType><xs:sequence><xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"/></xs:sequence></xs:complexType></xs:element><xs:element name="TermInfo"><xs:complexType><xs:sequence><xs:element ref="pc:TermName" minOccurs="0"/><xs:element ref="pc:TermId" minOccurs="0"/></xs:sequence></xs:complexType></xs:element><xs:element name="TermName" type="xs:string"/><xs:element name="TermId" type="xs:string"/></xs:schema></ct:contentTypeSchema>'

$updated = $false
$parts = $d.CustomXMLParts
if ($parts -ne $null) {
    $candidates = $parts.SelectByNamespace("http://schemas.microsoft.com/office/2006/metadata/contentType")
    if ($candidates -ne $null -and $candidates.Count -ge 1) {
        $part = $candidates.Item(1)
        $part.XML = $newSchemaXml
        $updated = $true
    }
    if (-not $updated) {
        for ($i = 1; $i -le $parts.Count; $i++) {
            $p = $parts.Item($i)
            if ($p.NamespaceURI -eq "http://schemas.microsoft.com/office/2006/metadata/contentType") {
                $p.XML = $newSchemaXml
                $updated = $true
            }
        }
    }
}

Write-Output "CustomXMLParts schema update attempted: $updated"
